# Reorder the ingredient lists ("Materias primas") stored in column C for
# several products, matching the shared-string reordering in the diff.
# No other cell values, layout or formatting changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Torta
$ws.Range("C2").Value = "3.0-Leche (litros),2.0-Harina  (kg),1.0-Vainilla (ml),1.0-Huevos (unidad),"

# Row 4: Queque
$ws.Range("C4").Value = "5.0-Harina  (kg),1.0-Vainilla (ml),2.0-Huevos (unidad),"

# Row 6: Pie de Limon
$ws.Range("C6").Value = "2.0-Limon (unidad),4.0-Harina  (kg),5.0-Crema (litros),5.0-Huevos (unidad),"

# Row 7: Cupcake
$ws.Range("C7").Value = "0.2-Leche (litros),0.3-Harina  (kg),0.1-Vainilla (ml),2.0-Huevos (unidad),"
